# "Atualizando a tabela de classificacao"
#
# Renames four club names in the standings table:
#   B6  "Bragantino"    -> "Red Bull Bragantino"
#   B7  "Athletico-PR"  -> "Athletico - PR"
#   B13 "Atlético-MG"   -> "Atlético - MG"
#   B18 "América-MG"    -> "América - MG"
#
# Because "Red Bull Bragantino" is noticeably longer than the other club
# names in the column, the cell picks up its own look (vertical="top" /
# horizontal="general" alignment, plus a slightly different hairline
# border) instead of the shared "left aligned" style used by every other
# team cell in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple text edits -----------------------------------------------
$ws.Range("B7").Value  = "Athletico - PR"
$ws.Range("B13").Value = "Atlético - MG"
$ws.Range("B18").Value = "América - MG"

# --- B6 needs new text + a distinct style -----------------------------
$b6 = $ws.Range("B6")
$b6.Value = "Red Bull Bragantino"

# Alignment: vertical top, general (not left) horizontal alignment.
$b6.VerticalAlignment   = -4160   # xlTop
$b6.HorizontalAlignment = 1       # xlGeneral

# Border: thin hairlines on all four edges - light grey (#EEEEEE) on the
# left/right/bottom, a touch darker (#DDDDDD) on top.
$b6.Borders.LineStyle     = 0            # xlLineStyleNone: start from a clean slate
$b6.Borders.Color         = 15658734     # RGB(238,238,238) = #EEEEEE on every edge
$b6.Borders.Item(8).Color = 14540253     # RGB(221,221,221) = #DDDDDD on the top edge only (xlEdgeTop = 8)
